$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2..11) for columns D, M, N, O, P, R, S
$data = @(
    @{ Row = 2;  D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; R = "Región de La Araucanía"; S = 2400 },
    @{ Row = 3;  D = 44574; M = 200; N = 3000; O = 3000; P = 3000; R = "Región de La Araucanía"; S = 3000 },
    @{ Row = 4;  D = 44616; M = 200; N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 },
    @{ Row = 5;  D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó"; S = 5000 },
    @{ Row = 6;  D = 44323; M = 20;  N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 },
    @{ Row = 7;  D = 44592; M = 5;   N = 7500; O = 7500; P = 7500; R = "Región de La Araucanía"; S = 7500 },
    @{ Row = 8;  D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; R = "Región de O'Higgins"; S = 3000 },
    @{ Row = 9;  D = 44551; M = 120; N = 4500; O = 4500; P = 4500; R = "Región de O'Higgins"; S = 4500 },
    @{ Row = 10; D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; R = "Región de La Araucanía"; S = 2800 },
    @{ Row = 11; D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; R = "Región de La Araucanía"; S = 1800 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 13).Value = $entry.M
    $ws.Cells.Item($r, 14).Value = $entry.N
    $ws.Cells.Item($r, 15).Value = $entry.O
    $ws.Cells.Item($r, 16).Value = $entry.P
    $ws.Cells.Item($r, 18).Value = $entry.R
    $ws.Cells.Item($r, 19).Value = $entry.S
}
